# Uren Registratie gister en vandaag
#
# Week 5's attendance table (rows 27-31 = Maandag..Vrijdag, columns C..I =
# Rief, Zinedine, Robin, Carlo, Marc, Sam, Michiel) gets hours filled in for
# "gisteren" (Dinsdag, row 28) and "vandaag" (Woensdag, row 29): everyone
# gets 4 hours on Dinsdag and 6 hours on Woensdag.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C28:I28").Value = 4
$ws.Range("C29:I29").Value = 6

# "Totaal lesuren p/w" for week 5 (B32) is a manually-entered total (unlike
# C32:I32, which are SUM formulas that recalc on their own) - bump it to
# match the new day totals.
$ws.Range("B32").Value = 14

# Leave the view scrolled/selected roughly where the edit happened.
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("K31").Select()
